$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: copy the existing rows 241 and 242 down to new rows 243 and 244
# (their original content is preserved as-is in the new rows).

# Row 243 = old Row 241 (Inferno / Arica y Parinacota / 2021-08-24)
$ws.Range("A243").Value = 4
$ws.Range("B243").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C243").Value = "Los Lagos"
$ws.Range("D243").Value = 44432
$ws.Range("D243").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E243").Value = 10
$ws.Range("F243").Value = 100112021
$ws.Range("G243").Value = "Ají"
$ws.Range("H243").Value = "Inferno"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 120
$ws.Range("K243").Value = 45000
$ws.Range("L243").Value = 45000
$ws.Range("M243").Value = 45000
$ws.Range("N243").Value = "`$/caja 12 kilos"
$ws.Range("O243").Value = "Región de Arica y Parinacota"
$ws.Range("P243").Value = 3750
$ws.Range("Q243").Value = 12
$ws.Range("R243").Value = "Hortaliza"

# Row 244 = old Row 242 (Inferno / Arica y Parinacota / 2022-01-31)
$ws.Range("A244").Value = 4
$ws.Range("B244").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C244").Value = "Los Lagos"
$ws.Range("D244").Value = 44592
$ws.Range("D244").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E244").Value = 10
$ws.Range("F244").Value = 100112021
$ws.Range("G244").Value = "Ají"
$ws.Range("H244").Value = "Inferno"
$ws.Range("I244").Value = "Primera"
$ws.Range("J244").Value = 35
$ws.Range("K244").Value = 17000
$ws.Range("L244").Value = 17000
$ws.Range("M244").Value = 17000
$ws.Range("N244").Value = "`$/caja 12 kilos"
$ws.Range("O244").Value = "Región de Arica y Parinacota"
$ws.Range("P244").Value = 1417
$ws.Range("Q244").Value = 12
$ws.Range("R244").Value = "Hortaliza"

# --- Step 2: overwrite rows 241 and 242 with the new weekly records.

# Row 241 -> Cristal / Región del Maule, dated 2022-04-05
$ws.Range("D241").Value = 44656
$ws.Range("H241").Value = "Cristal"
$ws.Range("J241").Value = 70
$ws.Range("K241").Value = 20000
$ws.Range("L241").Value = 20000
$ws.Range("M241").Value = 20000
$ws.Range("N241").Value = "`$/saco 25 kilos"
$ws.Range("O241").Value = "Región del Maule"
$ws.Range("P241").Value = 800
$ws.Range("Q241").Value = 25

# Row 242 -> still Inferno / Arica y Parinacota, but dated 2022-04-05 with new volumes/prices
$ws.Range("D242").Value = 44656
$ws.Range("J242").Value = 90
$ws.Range("K242").Value = 22000
$ws.Range("L242").Value = 22000
$ws.Range("M242").Value = 22000
$ws.Range("P242").Value = 1833
